$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3560.25
$ws.Range("I138").Value = 3096.2144
$ws.Range("J138").Value = 3731.2104
$ws.Range("K138").Value = 9288.643199999999
$ws.Range("L138").Value = 11193.6312
$ws.Range("M138").Value = -4148.643199999999
$ws.Range("N138").Value = -21473.6312
$ws.Range("H141").Value = 2440.8462
$ws.Range("I141").Value = 2311.9167
$ws.Range("K141").Value = 6935.750100000001
$ws.Range("M141").Value = -1755.750100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 10199.3
$ws.Range("I5").Value = 246.28572
$ws.Range("K5").Value = 246.28572
$ws.Range("M5").Value = -134.28572
$ws.Range("H32").Value = 9206.706
$ws.Range("I32").Value = 7557.2383
$ws.Range("J32").Value = 29990
$ws.Range("K32").Value = 7557.2383
$ws.Range("L32").Value = 29990
$ws.Range("M32").Value = -7270.2383
$ws.Range("N32").Value = -30564
$ws.Range("H61").Value = 6391.1577
$ws.Range("I61").Value = 7071.6206
$ws.Range("J61").Value = 4198.5557
$ws.Range("K61").Value = 7071.6206
$ws.Range("L61").Value = 4198.5557
$ws.Range("M61").Value = -6859.6206
$ws.Range("N61").Value = -4622.5557
$ws.Range("H97").Value = 3051161.8
$ws.Range("I97").Value = 3051161.8
$ws.Range("K97").Value = 3051161.8
$ws.Range("M97").Value = -3050665.8
$ws.Range("H122").Value = 3677743
$ws.Range("I122").Value = 6581445.5
$ws.Range("K122").Value = 19744336.5
$ws.Range("M122").Value = -19741886.5
$ws.Range("H132").Value = 21151.547
$ws.Range("I132").Value = 4365.1177
$ws.Range("J132").Value = 51190.42
$ws.Range("K132").Value = 13095.3531
$ws.Range("L132").Value = 153571.26
$ws.Range("M132").Value = -10565.3531
$ws.Range("N132").Value = -158631.26
$ws.Range("H136").Value = 6391.1577
$ws.Range("I136").Value = 7071.6206
$ws.Range("J136").Value = 4198.5557
$ws.Range("K136").Value = 21214.8618
$ws.Range("L136").Value = 12595.6671
$ws.Range("M136").Value = -18664.8618
$ws.Range("N136").Value = -17695.6671

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 10199.3
$ws.Range("I4").Value = 246.28572
$ws.Range("K4").Value = 246.28572
$ws.Range("M4").Value = -131.28572
$ws.Range("H20").Value = 3428.889
$ws.Range("I20").Value = 2809.4443
$ws.Range("K20").Value = 2809.4443
$ws.Range("M20").Value = -2562.4443
$ws.Range("H22").Value = 65.75
$ws.Range("I22").Value = 70.333336
$ws.Range("J22").Value = 52
$ws.Range("K22").Value = 70.333336
$ws.Range("L22").Value = 52
$ws.Range("M22").Value = 102.666664
$ws.Range("N22").Value = -398
$ws.Range("H134").Value = 9329.052
$ws.Range("I134").Value = 7171.6763
$ws.Range("K134").Value = 21515.0289
$ws.Range("M134").Value = -18980.0289

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 26666.666
$ws.Range("I6").Value = 26666.666
$ws.Range("K6").Value = 26666.666
$ws.Range("M6").Value = -26553.666
$ws.Range("H7").Value = 166.44444
$ws.Range("J7").Value = 700
$ws.Range("L7").Value = 700
$ws.Range("N7").Value = -926
$ws.Range("H31").Value = 19615.777
$ws.Range("I31").Value = 5338.2
$ws.Range("K31").Value = 5338.2
$ws.Range("M31").Value = -5043.2
$ws.Range("H34").Value = 19615.777
$ws.Range("I34").Value = 5338.2
$ws.Range("K34").Value = 5338.2
$ws.Range("M34").Value = -5136.2
$ws.Range("H41").Value = 9999.5
$ws.Range("I41").Value = 9999.5
$ws.Range("K41").Value = 9999.5
$ws.Range("M41").Value = -9571.5
$ws.Range("H51").Value = 56099
$ws.Range("J51").Value = 56099
$ws.Range("L51").Value = 56099
$ws.Range("N51").Value = -57571
$ws.Range("H61").Value = 56099
$ws.Range("J61").Value = 56099
$ws.Range("L61").Value = 56099
$ws.Range("N61").Value = -56795
$ws.Range("H132").Value = 65446.8
$ws.Range("I132").Value = 51833.55
$ws.Range("J132").Value = 92673.3
$ws.Range("K132").Value = 155500.65
$ws.Range("L132").Value = 278019.9
$ws.Range("M132").Value = -152970.65
$ws.Range("N132").Value = -283079.9
$ws.Range("H134").Value = 8315.161
$ws.Range("I134").Value = 5241.346
$ws.Range("K134").Value = 15724.038
$ws.Range("M134").Value = -13189.038
$ws.Range("H141").Value = 68757.37
$ws.Range("J141").Value = 71516.19500000001
$ws.Range("L141").Value = 71516.19500000001
$ws.Range("N141").Value = -81876.19500000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 45081.22
$ws.Range("I5").Value = 1051.9333
$ws.Range("J5").Value = 127636.125
$ws.Range("K5").Value = 3155.7999
$ws.Range("L5").Value = 382908.375
$ws.Range("M5").Value = -3043.7999
$ws.Range("N5").Value = -383132.375
$ws.Range("H23").Value = 727.8889
$ws.Range("I23").Value = 378.25
$ws.Range("K23").Value = 1134.75
$ws.Range("M23").Value = -899.75
$ws.Range("H26").Value = 132.4375
$ws.Range("I26").Value = 124.6
$ws.Range("K26").Value = 373.8
$ws.Range("M26").Value = -85.79999999999995
$ws.Range("H28").Value = 3000
$ws.Range("I28").Value = 3000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 9000
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("M28").Value = -8768
$ws.Range("H55").Value = 83336010
$ws.Range("I55").Value = 249751470
$ws.Range("J55").Value = 128275
$ws.Range("K55").Value = 749254410
$ws.Range("L55").Value = 384825
$ws.Range("M55").Value = -749254233
$ws.Range("N55").Value = -385179
$ws.Range("H59").Value = 4670
$ws.Range("I59").Value = 5026.25
$ws.Range("J59").Value = 3245
$ws.Range("K59").Value = 15078.75
$ws.Range("L59").Value = 9735
$ws.Range("M59").Value = -14538.75
$ws.Range("N59").Value = -10815
$ws.Range("H81").Value = 5376.6665
$ws.Range("I81").Value = 1319.25
$ws.Range("J81").Value = 6000.885
$ws.Range("K81").Value = 3957.75
$ws.Range("L81").Value = 18002.655
$ws.Range("M81").Value = -2834.75
$ws.Range("N81").Value = -20248.655
$ws.Range("H84").Value = 5376.6665
$ws.Range("I84").Value = 1319.25
$ws.Range("J84").Value = 6000.885
$ws.Range("K84").Value = 11873.25
$ws.Range("L84").Value = 54007.965
$ws.Range("M84").Value = -6257.25
$ws.Range("N84").Value = -65239.965
$ws.Range("H107").Value = 371.41666
$ws.Range("I107").Value = 327
$ws.Range("J107").Value = 403.14285
$ws.Range("K107").Value = 981
$ws.Range("L107").Value = 1209.42855
$ws.Range("M107").Value = 939
$ws.Range("N107").Value = -5049.428550000001
$ws.Range("H113").Value = 3075.121
$ws.Range("J113").Value = 1884.1
$ws.Range("L113").Value = 5652.299999999999
$ws.Range("N113").Value = -9992.299999999999
$ws.Range("H116").Value = 6879.8
$ws.Range("I116").Value = 821
$ws.Range("J116").Value = 8899.4
$ws.Range("K116").Value = 2463
$ws.Range("L116").Value = 26698.2
$ws.Range("M116").Value = 979
$ws.Range("N116").Value = -33582.2
$ws.Range("H131").Value = 14372597
$ws.Range("I131").Value = 27778446
$ws.Range("J131").Value = 12825768
$ws.Range("K131").Value = 83335338
$ws.Range("L131").Value = 38477304
$ws.Range("M131").Value = -83330298
$ws.Range("N131").Value = -38487384
$ws.Range("H132").Value = 1712.3077
$ws.Range("I132").Value = 1432.4166
$ws.Range("J132").Value = 1952.2142
$ws.Range("K132").Value = 12891.7494
$ws.Range("L132").Value = 17569.9278
$ws.Range("M132").Value = -10361.7494
$ws.Range("N132").Value = -22629.9278
$ws.Range("H135").Value = 45081.22
$ws.Range("I135").Value = 1051.9333
$ws.Range("J135").Value = 127636.125
$ws.Range("K135").Value = 9467.3997
$ws.Range("L135").Value = 1148725.125
$ws.Range("M135").Value = -6932.3997
$ws.Range("N135").Value = -1153795.125
$ws.Range("H140").Value = 2425.6667
$ws.Range("I140").Value = 1188.5
$ws.Range("K140").Value = 3565.5
$ws.Range("M140").Value = 1614.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("I70").Value = 20004296
$ws.Range("J70").Value = 26161.6
$ws.Range("K70").Value = 20004296
$ws.Range("L70").Value = 26161.6
$ws.Range("M70").Value = -20004026
$ws.Range("N70").Value = -26701.6
$ws.Range("I73").Value = 20004296
$ws.Range("J73").Value = 26161.6
$ws.Range("K73").Value = 20004296
$ws.Range("L73").Value = 26161.6
$ws.Range("M73").Value = -20003360
$ws.Range("N73").Value = -28033.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1377.1177
$ws.Range("I22").Value = 742.6
$ws.Range("K22").Value = 742.6
$ws.Range("M22").Value = -447.6
$ws.Range("H27").Value = 1377.1177
$ws.Range("I27").Value = 742.6
$ws.Range("K27").Value = 742.6
$ws.Range("M27").Value = -635.6
$ws.Range("H132").Value = 25471
$ws.Range("I132").Value = 28505.77
$ws.Range("K132").Value = 85517.31
$ws.Range("M132").Value = -82987.31
$ws.Range("H136").Value = 57050.133
$ws.Range("I136").Value = 70163.5
$ws.Range("K136").Value = 210490.5
$ws.Range("M136").Value = -207940.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1502
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 2004
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 2004
$ws.Range("M8").Value = -860
$ws.Range("N8").Value = -2284
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 500
$ws.Range("K17").Value = 500
$ws.Range("M17").Value = -328
$ws.Range("H107").Value = 41668308
$ws.Range("I107").Value = 100000860
$ws.Range("J107").Value = 2199
$ws.Range("K107").Value = 300002580
$ws.Range("L107").Value = 6597
$ws.Range("M107").Value = -300000660
$ws.Range("N107").Value = -10437
$ws.Range("H132").Value = 16312648
$ws.Range("I132").Value = 19238120
$ws.Range("J132").Value = 1100196.1
$ws.Range("K132").Value = 57714360
$ws.Range("L132").Value = 3300588.3
$ws.Range("M132").Value = -57711830
$ws.Range("N132").Value = -3305648.3
